$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.378.84"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3 (Ethereum)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.282.29"
$ws.Range("E3").Value = "  +1.50%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.21%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.28"
$ws.Range("E5").Value = "  +1.88%  "

# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.73"
$ws.Range("E6").Value = "  +2.28%  "

# Row 7 (USDC)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 (XRP)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9 (Dogecoin)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0957"
$ws.Range("E9").Value = "  +0.84%  "

# Row 10 (TRON)
$ws.Range("E10").Value = "  +0.96%  "

# Row 11 (Cardano)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("E11").Value = "  +3.99%  "

# Row 12 (Toncoin)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.90"
$ws.Range("E12").Value = "  +5.81%  "

# Row 15 (WrappedBTC)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.444.79"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16 (ShibaInu)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000130"
$ws.Range("E16").Value = "  +1.91%  "

# Row 17 (WrappedEther)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.291.86"
$ws.Range("E17").Value = "  +1.69%  "

# Row 18 (Chainlink)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.28"
$ws.Range("E18").Value = "  +3.23%  "

# Row 19 (Polkadot)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.15"
$ws.Range("E19").Value = "  +2.24%  "

# Row 20 (BitcoinCash)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "305.46"
$ws.Range("E20").Value = "  +2.47%  "

# Row 21 (Uniswap)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  -1.30%  "

# Row 22 (Dai)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23 (Litecoin)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.51"
$ws.Range("E23").Value = "  -2.06%  "

# Row 24 (Binance-PegBSC-USD)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.995"
$ws.Range("E24").Value = "  -2.00%  "

# Row 25 (Kaspa)
$ws.Range("E25").Value = "  +1.95%  "

# Row 26 (InternetComputer)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.37"
$ws.Range("E26").Value = "  +5.19%  "

# Row 27 (Monero)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.63"
$ws.Range("E27").Value = "  +4.10%  "

# Row 28 (PEPE)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0706"
$ws.Range("E28").Value = "  +5.00%  "

# Row 29 (Aptos)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.02"
$ws.Range("E29").Value = "  +3.63%  "

# Row 30 (PancakeSwap)
$ws.Range("E30").Value = "  +1.65%  "

# Row 31 (Fetch.AI)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("E31").Value = "  +3.80%  "

# Row 32 (USDe)
$ws.Range("E32").Value = "  +0.03%  "

# Row 33 (EthereumClassic)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.92"
$ws.Range("E33").Value = "  +1.97%  "

# Row 34 (FirstDigitalUSD)
$ws.Range("E34").Value = "  +0.04%  "

# Row 35 (SuiNetwork)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.930"
$ws.Range("E35").Value = "  +5.19%  "

# Row 36 (ImmutableX)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("E36").Value = "  +2.39%  "

# Row 37 (NEARProtocol)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  +3.18%  "

# Row 38 (PolygonEcosystemToken)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  +1.85%  "

# Row 39 (Stacks)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.41"
$ws.Range("E39").Value = "  +2.17%  "

# Row 40 (RenderToken)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.05"
$ws.Range("E40").Value = "  +3.45%  "

# Row 41 (Filecoin)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("E41").Value = "  +1.90%  "

# Row 42 (Aave)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.90"
$ws.Range("E42").Value = "  +0.24%  "

# Row 43 (Hedera)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0492"
$ws.Range("E43").Value = "  +2.76%  "

# Row 46 (Mantle)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.547"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47 (Polygon)
$ws.Range("E47").Value = "  +1.86%  "

# Row 48 (VeChain)
$ws.Range("E48").Value = "  +3.16%  "

# Row 49 (WhiteBITCoin)
$ws.Range("E49").Value = "  +0.69%  "

# Row 50 (InjectiveProtocol)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.17"
$ws.Range("E50").Value = "  +1.94%  "

# Row 51 (dogwifhat)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.56"
$ws.Range("E51").Value = "  +6.12%  "

# Row 13/14 swap: Avalanche <-> WrappedliquidstakedEther2.0
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.691.88"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.05"
$ws.Range("E14").Value = "  +6.50%  "

# Row 44/45 swap: Bittensor <-> Stellar
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0896"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "246.08"
$ws.Range("E45").Value = "  +5.05%  "
